# Revert "Actualizado el to-do"
# This reverts commit c825b7cd35743795125ff018770b0320811c7165.
#
# The previous commit had flipped three "Requerimientos" status cells from
# "p" to "x" and left the "Requerimientos" sheet/selection active; this
# reverts those cell values back to "p" and restores the "Must Have" sheet
# (with its prior selection) as the active sheet/selection.

$wb = $excel.ActiveWorkbook

# --- Restore the three status values on "Requerimientos" ---
$ws1 = $wb.Worksheets.Item("Requerimientos")
$ws1.Range("B2").Value = "p"
$ws1.Range("B10").Value = "p"
$ws1.Range("B11").Value = "p"

# Restore the previously-saved selection on "Requerimientos" (no longer the
# active sheet, but Excel still persists a remembered selection per sheet).
[void]$ws1.Range("B2").Select()

# --- Make "Must Have" the active sheet again, with its prior selection ---
$ws2 = $wb.Worksheets.Item("Must Have")
[void]$ws2.Activate()
[void]$ws2.Range("D9").Select()
